{"js": "const body = context.document.body;\n\nconst replacements = [\n  {\n    search: \"Metamodel messages: (match, apply) CSPO quads for each Resource hierarchy new instance: quads message. Apply occurrences to each local matching CSPO. Metamodels aggregate new occurrences.\",\n    replace: \"Metamodel messages: (match, apply) CSPO quads for each Resource hierarchy new instance: quads message. Apply occurrences to each local matching CSPO. Context of each applied CSPO: complement triple (i.e.: CPO for S) resources history. Metamodels aggregate new occurrences.\"\n  },\n  {\n    search: \"Resource history: invoked (match, apply) transforms in contexts until base resources.\",\n    replace: \"Resource history: invoked (match, apply) transforms in contexts until base resources. Complement based ID encoding.\"\n  }\n];\n\nfor (const { search, replace } of replacements) {\n  const results = body.search(search, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + search);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Metamodel messages: (match, apply) CSPO quads for each Resource hierarchy new instance: quads message. Apply occurrences to each local matching CSPO. Metamodels aggregate new occurrences.\"\n$find.Replacement.Text = \"Metamodel messages: (match, apply) CSPO quads for each Resource hierarchy new instance: quads message. Apply occurrences to each local matching CSPO. Context of each applied CSPO: complement triple (i.e.: CPO for S) resources history. Metamodels aggregate new occurrences.\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Resource history: invoked (match, apply) transforms in contexts until base resources.\"\n$find2.Replacement.Text = \"Resource history: invoked (match, apply) transforms in contexts until base resources. Complement based ID encoding.\"\n$find2.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n"}
